$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Preserve the cell as text (matches the source data's inlineStr/string
    # typing) instead of letting Excel auto-coerce the numeric-looking
    # string into a Number cell. NumberFormat="@" forces text entry, and
    # Style="Normal" resets back to the default (unstyled) cell format so
    # no stray number-format style is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("C2") "445"
Set-TextValue $ws.Range("D2") "1051828.79"
Set-TextValue $ws.Range("C4") "910"
Set-TextValue $ws.Range("D4") "3018233.54"
Set-TextValue $ws.Range("C6") "564"
Set-TextValue $ws.Range("D6") "1705700.55"
Set-TextValue $ws.Range("C7") "14"
Set-TextValue $ws.Range("D7") "31000.00"
Set-TextValue $ws.Range("C8") "29"
Set-TextValue $ws.Range("D8") "116938.45"
Set-TextValue $ws.Range("C14") "212"
Set-TextValue $ws.Range("D14") "566362.00"
Set-TextValue $ws.Range("C16") "478"
Set-TextValue $ws.Range("D16") "1729608.75"
Set-TextValue $ws.Range("C17") "138"
Set-TextValue $ws.Range("D17") "420415.33"
Set-TextValue $ws.Range("C20") "170"
Set-TextValue $ws.Range("D20") "430599.00"
Set-TextValue $ws.Range("C21") "331"
Set-TextValue $ws.Range("D21") "1167741.00"
Set-TextValue $ws.Range("C22") "157"
Set-TextValue $ws.Range("D22") "444137.26"
Set-TextValue $ws.Range("C23") "9"
Set-TextValue $ws.Range("D23") "21000.00"
Set-TextValue $ws.Range("C24") "11"
Set-TextValue $ws.Range("D24") "48900.00"
Set-TextValue $ws.Range("C28") "263"
Set-TextValue $ws.Range("D28") "666542.64"
Set-TextValue $ws.Range("C30") "520"
Set-TextValue $ws.Range("D30") "2057800.70"
Set-TextValue $ws.Range("C32") "369"
Set-TextValue $ws.Range("D32") "1215198.17"
Set-TextValue $ws.Range("C33") "13"
Set-TextValue $ws.Range("D33") "41500.00"
Set-TextValue $ws.Range("C34") "15"
Set-TextValue $ws.Range("D34") "45932.00"
Set-TextValue $ws.Range("C40") "138"
Set-TextValue $ws.Range("D40") "387482.22"
Set-TextValue $ws.Range("C41") "82"
Set-TextValue $ws.Range("D41") "397909.98"
Set-TextValue $ws.Range("C42") "121"
Set-TextValue $ws.Range("D42") "502119.99"
Set-TextValue $ws.Range("C45") "348"
Set-TextValue $ws.Range("D45") "958867.74"
Set-TextValue $ws.Range("C47") "583"
Set-TextValue $ws.Range("D47") "2255408.99"
Set-TextValue $ws.Range("C48") "394"
Set-TextValue $ws.Range("D48") "1310867.16"
Set-TextValue $ws.Range("C51") "3554"
Set-TextValue $ws.Range("D51") "8111876.70"
Set-TextValue $ws.Range("C53") "3865"
Set-TextValue $ws.Range("D53") "13217200.06"
Set-TextValue $ws.Range("C55") "3952"
Set-TextValue $ws.Range("D55") "12149356.47"
Set-TextValue $ws.Range("C56") "55"
Set-TextValue $ws.Range("D56") "145350.00"
Set-TextValue $ws.Range("C73") "377"
Set-TextValue $ws.Range("D73") "934971.70"
Set-TextValue $ws.Range("C74") "6"
Set-TextValue $ws.Range("D74") "32254.00"
Set-TextValue $ws.Range("C75") "906"
Set-TextValue $ws.Range("D75") "3081919.39"
Set-TextValue $ws.Range("C76") "510"
Set-TextValue $ws.Range("D76") "1662286.87"
Set-TextValue $ws.Range("C77") "37"
Set-TextValue $ws.Range("D77") "102500.00"
